$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the bottom of the table (8 and 9), which in Excel
# naturally inherit the formatting of the row above (row 7).
$ws.Rows("8:9").Insert() | Out-Null

# Row 8: GSM8K / gsm8-std / gemini-2.5-flash-preview-05-20
# (benchmark_name "gsm8-std" is registered as a shared string before the
# dataset name "GSM8K", matching the original authoring order.)
$ws.Range("B8").Value = "gsm8-std"
$ws.Range("A8").Value = "GSM8K"
$ws.Range("C8").Value = "gemini-2.5-flash-preview-05-20"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 96
$ws.Range("F8").Value = 100
$ws.Range("G8").Value = 98
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 97
$ws.Range("L8").Value = 98
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0.96
$ws.Range("P8").Value = 0.98
$ws.Range("Q8").Value = 0.01
$ws.Range("R8").Value = 0.01
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0.97
$ws.Range("U8").Value = 0.98

# Row 9: GSM8K / gsm8-std / gemini-2.5-pro-preview-05-06
$ws.Range("A9").Value = "GSM8K"
$ws.Range("B9").Value = "gsm8-std"
$ws.Range("C9").Value = "gemini-2.5-pro-preview-05-06"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 96
$ws.Range("F9").Value = 100
$ws.Range("G9").Value = 97
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 98
$ws.Range("L9").Value = 99
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0.96
$ws.Range("P9").Value = 0.97
$ws.Range("Q9").Value = 0.02
$ws.Range("R9").Value = 0.01
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0.98
$ws.Range("U9").Value = 0.99

# Rows 6 and 7 lose their explicit (redundant) cell formatting.
$ws.Range("A6:U7").ClearFormats()

# Leave the selection where the user finished editing.
$ws.Range("D9").Select() | Out-Null
